# Update the Norfolk St._B team-specific transition-probability matrix
# with recalculated values after adding team-specific time data
# (logic for using this data in simulation not yet implemented).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2049689440993789
$ws.Range("C2").Value = 0.5403726708074534
$ws.Range("J2").Value = 0.009316770186335404
$ws.Range("P2").Value = 0.1583850931677019
$ws.Range("S2").Value = 0.08695652173913043
$ws.Range("B3").Value = 0.01092896174863388
$ws.Range("C3").Value = 0.04918032786885246
$ws.Range("J3").Value = 0.0546448087431694
$ws.Range("P3").Value = 0.7540983606557377
$ws.Range("S3").Value = 0.1311475409836066
$ws.Range("J4").Value = 0.03333333333333333
$ws.Range("O4").Value = 0.03333333333333333
$ws.Range("P4").Value = 0.6333333333333333
$ws.Range("S4").Value = 0.3
$ws.Range("B6").Value = 0.0861244019138756
$ws.Range("D6").Value = 0.004784688995215311
$ws.Range("E6").Value = 0.004784688995215311
$ws.Range("F6").Value = 0.06220095693779904
$ws.Range("J6").Value = 0.3062200956937799
$ws.Range("O6").Value = 0.02870813397129187
$ws.Range("Q6").Value = 0.07655502392344497
$ws.Range("R6").Value = 0.07655502392344497
$ws.Range("S6").Value = 0.354066985645933
$ws.Range("B7").Value = 0.13
$ws.Range("D7").Value = 0.02
$ws.Range("E7").Value = 0.005
$ws.Range("F7").Value = 0.065
$ws.Range("J7").Value = 0.155
$ws.Range("O7").Value = 0.02
$ws.Range("Q7").Value = 0.13
$ws.Range("R7").Value = 0.07000000000000001
$ws.Range("S7").Value = 0.405
$ws.Range("B8").Value = 0.1333333333333333
$ws.Range("D8").Value = 0.01428571428571429
$ws.Range("F8").Value = 0.07857142857142857
$ws.Range("J8").Value = 0.1023809523809524
$ws.Range("O8").Value = 0.01666666666666667
$ws.Range("Q8").Value = 0.1285714285714286
$ws.Range("R8").Value = 0.1142857142857143
$ws.Range("S8").Value = 0.4119047619047619
$ws.Range("B9").Value = 0.08695652173913043
$ws.Range("D9").Value = 0.01449275362318841
$ws.Range("F9").Value = 0.05797101449275362
$ws.Range("J9").Value = 0.07246376811594203
$ws.Range("O9").Value = 0.02898550724637681
$ws.Range("Q9").Value = 0.1594202898550725
$ws.Range("R9").Value = 0.09420289855072464
$ws.Range("S9").Value = 0.4855072463768116
$ws.Range("B10").Value = 0.1286113699906803
$ws.Range("D10").Value = 0.02143522833178006
$ws.Range("E10").Value = 0.001863932898415657
$ws.Range("F10").Value = 0.05964585274930102
$ws.Range("J10").Value = 0.1164958061509786
$ws.Range("O10").Value = 0.01211556383970177
$ws.Range("Q10").Value = 0.1602982292637465
$ws.Range("R10").Value = 0.07921714818266543
$ws.Range("S10").Value = 0.4203168685927307
$ws.Range("G11").Value = 0.1301587301587302
$ws.Range("J11").Value = 0.08888888888888889
$ws.Range("K11").Value = 0.1841269841269841
$ws.Range("L11").Value = 0.5777777777777777
$ws.Range("S11").Value = 0.01904761904761905
$ws.Range("G12").Value = 0.7431693989071039
$ws.Range("J12").Value = 0.1912568306010929
$ws.Range("K12").Value = 0.00546448087431694
$ws.Range("L12").Value = 0.01092896174863388
$ws.Range("S12").Value = 0.04918032786885246
$ws.Range("G13").Value = 0.7105263157894737
$ws.Range("J13").Value = 0.2631578947368421
$ws.Range("S13").Value = 0.02631578947368421
$ws.Range("F15").Value = 0.0273224043715847
$ws.Range("H15").Value = 0.1475409836065574
$ws.Range("I15").Value = 0.0546448087431694
$ws.Range("J15").Value = 0.4426229508196721
$ws.Range("K15").Value = 0.04918032786885246
$ws.Range("M15").Value = 0.00546448087431694
$ws.Range("O15").Value = 0.0546448087431694
$ws.Range("S15").Value = 0.2185792349726776
$ws.Range("F16").Value = 0.04368932038834952
$ws.Range("H16").Value = 0.1067961165048544
$ws.Range("I16").Value = 0.05825242718446602
$ws.Range("J16").Value = 0.4223300970873786
$ws.Range("K16").Value = 0.1650485436893204
$ws.Range("M16").Value = 0.01456310679611651
$ws.Range("O16").Value = 0.05825242718446602
$ws.Range("S16").Value = 0.1310679611650485
$ws.Range("F17").Value = 0.04166666666666666
$ws.Range("H17").Value = 0.1597222222222222
$ws.Range("I17").Value = 0.08333333333333333
$ws.Range("J17").Value = 0.3680555555555556
$ws.Range("K17").Value = 0.1493055555555556
$ws.Range("M17").Value = 0.01041666666666667
$ws.Range("O17").Value = 0.05208333333333334
$ws.Range("S17").Value = 0.1354166666666667
$ws.Range("F18").Value = 0.02777777777777778
$ws.Range("H18").Value = 0.2333333333333333
$ws.Range("I18").Value = 0.1111111111111111
$ws.Range("J18").Value = 0.3666666666666666
$ws.Range("K18").Value = 0.06666666666666667
$ws.Range("M18").Value = 0.02777777777777778
$ws.Range("N18").Value = 0.005555555555555556
$ws.Range("O18").Value = 0.05555555555555555
$ws.Range("S18").Value = 0.1055555555555556
$ws.Range("F19").Value = 0.02005012531328321
$ws.Range("H19").Value = 0.2330827067669173
$ws.Range("I19").Value = 0.06265664160401002
$ws.Range("J19").Value = 0.3316624895572264
$ws.Range("K19").Value = 0.1261487050960735
$ws.Range("M19").Value = 0.02255639097744361
$ws.Range("N19").Value = 0.001670843776106934
$ws.Range("O19").Value = 0.06516290726817042
$ws.Range("S19").Value = 0.1370091896407686
